# Weekly update: insert a new daily price record for
# "Vega Monumental Concepción - Cilantro" (Hortaliza) ahead of the
# existing rows, pushing the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new blank row at row 333; this shifts the previous
# rows 333..359 down to 334..360 (and grows the sheet dimension to
# A1:R360), exactly like Excel's native "Insert Row" behaviour.
$ws.Rows(333).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A333").Value = 11
$ws.Range("B333").Value = "Vega Monumental Concepción"
$ws.Range("C333").Value = "Bíobío"
$ws.Range("D333").Value = 45223
$ws.Range("E333").Value = 8
$ws.Range("F333").Value = 100112040
$ws.Range("G333").Value = "Cilantro"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 60
$ws.Range("K333").Value = 10000
$ws.Range("L333").Value = 10000
$ws.Range("M333").Value = 10000
$ws.Range("N333").Value = "$/caja 36 atados"
$ws.Range("O333").Value = "Región Metropolitana"
$ws.Range("P333").Value = 278
$ws.Range("Q333").Value = 36
$ws.Range("R333").Value = "Hortaliza"
